# chore: adapt column header formatting to respective input file names
#
# 1. Rename the "_old"/"_new" suffixes on the header row (row 1) to
#    "_FV2304"/"_FV2310" respectively.
# 2. Turn the used range (A1:U58) into a real Excel Table (ListObject)
#    so the renamed headers become the table's column names, and an
#    autofilter is shown on row 1.
# 3. Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 21   # column U
$headerRow = 1

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item($headerRow, $c)
    $v = $cell.Value2
    if ($v -ne $null) {
        $newValue = $v
        if ($v.EndsWith("_old")) {
            $newValue = $v.Substring(0, $v.Length - 4) + "_FV2304"
        } elseif ($v.EndsWith("_new")) {
            $newValue = $v.Substring(0, $v.Length - 4) + "_FV2310"
        }
        if ($newValue -ne $v) {
            $cell.Value2 = $newValue
        }
    }
}

# Turn the data range into a native Excel table, picking up the
# (now renamed) header row as the column headers, with an autofilter.
$tableRange = $ws.Range("A1:U58")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.TableStyle = "TableStyleNone"

# Freeze the header row (row 1) via the top-left cell of the scrollable area.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
